# Auto-generated edit script: update crypto price (D) and 1h volume change (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.492.18"
$ws.Range("E2").Value = "  -1.26%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.763.06"
$ws.Range("E3").Value = "  -2.48%  "

$ws.Range("E4").Value = "  +0.22%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.22"
$ws.Range("E5").Value = "  -0.93%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.32"
$ws.Range("E6").Value = "  +0.45%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.759.85"
$ws.Range("E7").Value = "  -2.56%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("E9").Value = "  -0.17%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.164"
$ws.Range("E10").Value = "  +0.26%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.48"
$ws.Range("E11").Value = "  +0.28%  "

$ws.Range("E12").Value = "  -0.23%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000271"
$ws.Range("E13").Value = "  +3.76%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.79"
$ws.Range("E14").Value = "  -0.74%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.396.06"
$ws.Range("E15").Value = "  -2.35%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.796.61"
$ws.Range("E16").Value = "  -1.48%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.99"
$ws.Range("E17").Value = "  +4.89%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.482.26"
$ws.Range("E18").Value = "  -1.25%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.24"
$ws.Range("E19").Value = "  -1.79%  "

$ws.Range("E20").Value = "  +0.90%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.53"
$ws.Range("E21").Value = "  -2.55%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "468.46"
$ws.Range("E22").Value = "  +0.49%  "

$ws.Range("E23").Value = "  -1.48%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000149"
$ws.Range("E24").Value = "  -6.69%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.92"
$ws.Range("E25").Value = "  +1.21%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.23"
$ws.Range("E26").Value = "  -0.17%  "

$ws.Range("E27").Value = "  +0.83%  "

$ws.Range("E28").Value = "  +2.63%  "

$ws.Range("E29").Value = "  -0.02%  "

$ws.Range("E30").Value = "  -2.20%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.918.85"
$ws.Range("E31").Value = "  -2.15%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.62"
$ws.Range("E32").Value = "  +0.04%  "

$ws.Range("E33").Value = "  -2.89%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.40"
$ws.Range("E34").Value = "  -2.78%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.17"
$ws.Range("E35").Value = "  -3.84%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.729.08"
$ws.Range("E36").Value = "  -2.36%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.83"
$ws.Range("E37").Value = "  +4.51%  "

$ws.Range("E38").Value = "  +0.41%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.89"
$ws.Range("E39").Value = "  -0.36%  "

$ws.Range("E40").Value = "  -2.28%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.995"
$ws.Range("E41").Value = "  -2.79%  "

$ws.Range("E42").Value = "  +0.08%  "

$ws.Range("E43").Value = "  +0.14%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.72"
$ws.Range("E45").Value = "  +1.21%  "

$ws.Range("E46").Value = "  -1.01%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "46.24"
$ws.Range("E47").Value = "  -1.86%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "401.50"
$ws.Range("E48").Value = "  -4.75%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000273"
$ws.Range("E49").Value = "  -7.64%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "141.81"
$ws.Range("E50").Value = "  -0.74%  "

$ws.Range("E51").Value = "  -1.14%  "
